$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while forcing text storage so that
# numeric-looking strings (prices, percentages) are not auto-converted
# to numbers/dates by Excel, and restore the default style afterwards
# so no stray number-format style is left on the cell.
function Set-TextCell {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Rows 2-37: Price/Volume(1h) updates ---
Set-TextCell "D2" '37.054.79'
Set-TextCell "E2" '  +0.34%  '
Set-TextCell "D3" '2.052.22'
Set-TextCell "E3" '  +4.18%  '
Set-TextCell "E4" '  -0.15%  '
Set-TextCell "D5" '234.69'
Set-TextCell "E5" '  -1.46%  '
Set-TextCell "D6" '0.608'
Set-TextCell "E6" '  +2.12%  '
Set-TextCell "E7" '  -0.04%  '
Set-TextCell "D8" '56.17'
Set-TextCell "E8" '  +4.99%  '
Set-TextCell "D9" '0.377'
Set-TextCell "E9" '  +2.73%  '
Set-TextCell "D10" '57.90'
Set-TextCell "E10" '  -2.37%  '
Set-TextCell "D11" '0.0760'
Set-TextCell "E11" '  +2.36%  '
Set-TextCell "E12" '  +3.12%  '
Set-TextCell "D13" '2.349.24'
Set-TextCell "E13" '  +4.02%  '
Set-TextCell "D14" '14.46'
Set-TextCell "E14" '  +5.09%  '
Set-TextCell "D15" '20.54'
Set-TextCell "E15" '  -0.59%  '
Set-TextCell "D16" '0.771'
Set-TextCell "E16" '  +3.70%  '
Set-TextCell "D17" '5.21'
Set-TextCell "E17" '  +4.32%  '
Set-TextCell "D18" '2.043.21'
Set-TextCell "E18" '  +3.90%  '
Set-TextCell "D19" '37.101.51'
Set-TextCell "E19" '  +0.81%  '
Set-TextCell "D20" '5.86'
Set-TextCell "E20" '  +19.79%  '
Set-TextCell "D21" '68.16'
Set-TextCell "E21" '  +0.89%  '
Set-TextCell "D22" '0.0₃0806'
Set-TextCell "E22" '  +0.54%  '
Set-TextCell "D23" '222.29'
Set-TextCell "E23" '  -1.61%  '
Set-TextCell "E24" '  -0.15%  '
Set-TextCell "D25" '2.41'
Set-TextCell "E25" '  +4.42%  '
Set-TextCell "D26" '2.42'
Set-TextCell "E26" '  +2.44%  '
Set-TextCell "D27" '163.17'
Set-TextCell "E27" '  +1.24%  '
Set-TextCell "D28" '8.82'
Set-TextCell "E28" '  +3.40%  '
Set-TextCell "D29" '0.129'
Set-TextCell "E29" '  +2.91%  '
Set-TextCell "D30" '19.17'
Set-TextCell "E30" '  +1.28%  '
Set-TextCell "D31" '1.36'
Set-TextCell "E31" '  +8.05%  '
Set-TextCell "E32" '  +1.15%  '
Set-TextCell "D33" '4.42'
Set-TextCell "E33" '  +1.72%  '
Set-TextCell "D34" '0.0613'
Set-TextCell "E34" '  +1.45%  '
Set-TextCell "D35" '2.49'
Set-TextCell "E35" '  +8.66%  '
Set-TextCell "D36" '4.33'
Set-TextCell "E36" '  +2.74%  '
Set-TextCell "E37" '  -0.19%  '

# --- Rows 38-40: coin order rotated (THORChain, RenderToken, WEMIXToken) ---
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell "D38" '5.89'
Set-TextCell "E38" '  +16.20%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell "D39" '3.28'
Set-TextCell "E39" '  +1.30%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell "D40" '1.75'
Set-TextCell "E40" '  -1.43%  '

# --- Rows 41-42: Price/Volume(1h) updates ---
Set-TextCell "E41" '  -2.46%  '
Set-TextCell "D42" '4.46'
Set-TextCell "E42" '  +25.59%  '

# --- Rows 43-44: coin order swapped (Cronos, Maker) ---
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell "D43" '0.0953'
Set-TextCell "E43" '  +7.72%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell "D44" '1.483.44'
Set-TextCell "E44" '  +6.19%  '

# --- Rows 45-51: Price/Volume(1h) updates ---
Set-TextCell "D45" '94.38'
Set-TextCell "E45" '  +9.47%  '
Set-TextCell "D46" '0.0207'
Set-TextCell "E46" '  +2.58%  '
Set-TextCell "D47" '1.13'
Set-TextCell "E47" '  +0.81%  '
Set-TextCell "D48" '15.89'
Set-TextCell "E48" '  +6.02%  '
Set-TextCell "D49" '1.02'
Set-TextCell "E49" '  +3.09%  '
Set-TextCell "D50" '2.91'
Set-TextCell "E50" '  +2.03%  '
Set-TextCell "D51" '7.02'
Set-TextCell "E51" '  +6.88%  '
